$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 data rows (rows 4-6), shifting nothing up since they are the tail
$ws.Range("A4:E6").Delete()

# Update remaining two data rows with the reduced "10V" recipe/morbidity set
$ws.Cells.Item(2, 1).Value = "10VB1AS1BloodTest    "
$ws.Cells.Item(2, 3).Value = "10VA1C1spondylosis       "

$ws.Cells.Item(3, 1).Value = "10VB2AS1BloodTest    "
$ws.Cells.Item(3, 3).Value = "10VA2C1spondylosis       "

# Blank out column E for the remaining rows, but keep the cells present (empty)
$ws.Range("E2:E3").ClearContents()
$ws.Range("E2:E3").Font.Bold = $false

# Widen column A to fit the longer strings
$ws.Columns.Item(1).ColumnWidth = 23

# Move the selection to B2
$ws.Range("B2").Select()
